# Proveedores - Orden compra Excel
# Adds "PROVEEDOR ÚLTIMA COMPRA" column to the "Movimientos" sheet and
# renames the "Tipo docto (FOC)" header (column B) to simply "Tipo docto".

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Movimientos")

# --- header text edits -----------------------------------------------
# Write the brand-new column first so its shared string is appended
# before the renamed "Tipo docto" string (keeps shared-string order
# stable / matches how Excel appends new unique strings as they are
# first encountered).
$ws.Range("P1").Value = "PROVEEDOR ÚLTIMA COMPRA"
$ws.Range("B1").Value = "Tipo docto"

# --- header row layout --------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30

$header = $ws.Range("A1:P1")
$header.VerticalAlignment = -4108   # xlCenter
$header.WrapText = $true

# --- column widths (characters) -----------------------------------------
# ColumnWidth is expressed in "characters"; Excel stores the sheet's raw
# <col width="..."/> as characters + 5/6 (the default-font padding), so
# each value below is pre-offset to land on the target stored width.
$ws.Columns.Item(1).ColumnWidth  = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth  = 10.736979166666666
$ws.Columns.Item(3).ColumnWidth  = 13.736979166666666
$ws.Columns.Item(4).ColumnWidth  = 16.592447916666668
$ws.Columns.Item(6).ColumnWidth  = 12.877604166666666
$ws.Columns.Item(9).ColumnWidth  = 9.736979166666666
$ws.Columns.Item(10).ColumnWidth = 12.451822916666666
$ws.Columns.Item(11).ColumnWidth = 21.166666666666668
$ws.Columns.Item(16).ColumnWidth = 26.877604166666668
